$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 20790
$ws.Cells.Item(62, 9).Value = 33533.332
$ws.Cells.Item(62, 10).Value = 1675
$ws.Cells.Item(62, 11).Value = 33533.332
$ws.Cells.Item(62, 12).Value = 1675
$ws.Cells.Item(62, 13).Value = -32909.332
$ws.Cells.Item(62, 14).Value = -2923

$ws.Cells.Item(65, 8).Value = 20790
$ws.Cells.Item(65, 9).Value = 33533.332
$ws.Cells.Item(65, 10).Value = 1675
$ws.Cells.Item(65, 11).Value = 167666.66
$ws.Cells.Item(65, 12).Value = 8375
$ws.Cells.Item(65, 13).Value = -164546.66
$ws.Cells.Item(65, 14).Value = -14615

$ws.Cells.Item(86, 8).Value = 5842.375
$ws.Cells.Item(86, 9).Value = 5769.5
$ws.Cells.Item(86, 10).Value = 5866.6665
$ws.Cells.Item(86, 11).Value = 5769.5
$ws.Cells.Item(86, 12).Value = 5866.6665
$ws.Cells.Item(86, 13).Value = -4646.5
$ws.Cells.Item(86, 14).Value = -8112.6665

$ws.Cells.Item(89, 8).Value = 5842.375
$ws.Cells.Item(89, 9).Value = 5769.5
$ws.Cells.Item(89, 10).Value = 5866.6665
$ws.Cells.Item(89, 11).Value = 28847.5
$ws.Cells.Item(89, 12).Value = 29333.3325
$ws.Cells.Item(89, 13).Value = -23231.5
$ws.Cells.Item(89, 14).Value = -40565.3325

$ws.Cells.Item(107, 8).Value = 1439.5555
$ws.Cells.Item(107, 9).Value = 1156.6471
$ws.Cells.Item(107, 10).Value = 1920.5
$ws.Cells.Item(107, 11).Value = 1156.6471
$ws.Cells.Item(107, 12).Value = 1920.5
$ws.Cells.Item(107, 13).Value = 763.3529000000001
$ws.Cells.Item(107, 14).Value = -5760.5

$ws.Cells.Item(111, 8).Value = 986.56525
$ws.Cells.Item(111, 9).Value = 546.6667
$ws.Cells.Item(111, 10).Value = 1466.4546
$ws.Cells.Item(111, 11).Value = 1640.0001
$ws.Cells.Item(111, 12).Value = 4399.3638
$ws.Cells.Item(111, 13).Value = 1426.9999
$ws.Cells.Item(111, 14).Value = -10533.3638

$ws.Cells.Item(129, 8).Value = 894.0925999999999
$ws.Cells.Item(129, 9).Value = 600
$ws.Cells.Item(129, 11).Value = 1800
$ws.Cells.Item(129, 13).Value = 3200

$ws.Cells.Item(138, 8).Value = 2261.0203
$ws.Cells.Item(138, 9).Value = 1798.5
$ws.Cells.Item(138, 10).Value = 2270.5566
$ws.Cells.Item(138, 11).Value = 5395.5
$ws.Cells.Item(138, 12).Value = 6811.6698
$ws.Cells.Item(138, 13).Value = -255.5
$ws.Cells.Item(138, 14).Value = -17091.6698


$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1802.5641
$ws.Cells.Item(32, 9).Value = 1864.6216
$ws.Cells.Item(32, 10).Value = 654.5
$ws.Cells.Item(32, 11).Value = 1864.6216
$ws.Cells.Item(32, 12).Value = 654.5
$ws.Cells.Item(32, 13).Value = -1577.6216
$ws.Cells.Item(32, 14).Value = -1228.5

$ws.Cells.Item(74, 8).Value = 773.9394
$ws.Cells.Item(74, 9).Value = 721.3570999999999
$ws.Cells.Item(74, 10).Value = 1068.4
$ws.Cells.Item(74, 11).Value = 721.3570999999999
$ws.Cells.Item(74, 12).Value = 1068.4
$ws.Cells.Item(74, 13).Value = 152.6429000000001
$ws.Cells.Item(74, 14).Value = -2816.4

$ws.Cells.Item(77, 8).Value = 773.9394
$ws.Cells.Item(77, 9).Value = 721.3570999999999
$ws.Cells.Item(77, 10).Value = 1068.4
$ws.Cells.Item(77, 11).Value = 3606.7855
$ws.Cells.Item(77, 12).Value = 5342
$ws.Cells.Item(77, 13).Value = 761.2145
$ws.Cells.Item(77, 14).Value = -14078


$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 125002660
$ws.Cells.Item(105, 9).Value = 166669040
$ws.Cells.Item(105, 11).Value = 166669040
$ws.Cells.Item(105, 13).Value = -166667293


$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 866.65576
$ws.Cells.Item(31, 9).Value = 764.125
$ws.Cells.Item(31, 10).Value = 1061.9524
$ws.Cells.Item(31, 11).Value = 764.125
$ws.Cells.Item(31, 12).Value = 1061.9524
$ws.Cells.Item(31, 13).Value = -469.125
$ws.Cells.Item(31, 14).Value = -1651.9524

$ws.Cells.Item(34, 8).Value = 866.65576
$ws.Cells.Item(34, 9).Value = 764.125
$ws.Cells.Item(34, 10).Value = 1061.9524
$ws.Cells.Item(34, 11).Value = 764.125
$ws.Cells.Item(34, 12).Value = 1061.9524
$ws.Cells.Item(34, 13).Value = -562.125
$ws.Cells.Item(34, 14).Value = -1465.9524

$ws.Cells.Item(107, 8).Value = 679.6087
$ws.Cells.Item(107, 9).Value = 556.8
$ws.Cells.Item(107, 10).Value = 774.0769
$ws.Cells.Item(107, 11).Value = 556.8
$ws.Cells.Item(107, 12).Value = 774.0769
$ws.Cells.Item(107, 13).Value = 1363.2
$ws.Cells.Item(107, 14).Value = -4614.0769


$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(33, 8).Value = 161
$ws.Cells.Item(33, 9).Value = 111
$ws.Cells.Item(33, 10).Value = 231
$ws.Cells.Item(33, 11).Value = 666
$ws.Cells.Item(33, 12).Value = 1386
$ws.Cells.Item(33, 13).Value = -383
$ws.Cells.Item(33, 14).Value = -1952

$ws.Cells.Item(64, 8).Value = 3581.111
$ws.Cells.Item(64, 10).Value = 4538
$ws.Cells.Item(64, 12).Value = 13614
$ws.Cells.Item(64, 14).Value = -14154

$ws.Cells.Item(67, 8).Value = 3581.111
$ws.Cells.Item(67, 10).Value = 4538
$ws.Cells.Item(67, 12).Value = 13614
$ws.Cells.Item(67, 14).Value = -15486

$ws.Cells.Item(68, 8).Value = 1368.5483
$ws.Cells.Item(68, 10).Value = 1925.8823
$ws.Cells.Item(68, 12).Value = 5777.6469
$ws.Cells.Item(68, 14).Value = -7399.6469

$ws.Cells.Item(71, 8).Value = 1368.5483
$ws.Cells.Item(71, 10).Value = 1925.8823
$ws.Cells.Item(71, 12).Value = 17332.9407
$ws.Cells.Item(71, 14).Value = -25444.9407

$ws.Cells.Item(98, 8).Value = 804
$ws.Cells.Item(98, 9).Value = 116.75
$ws.Cells.Item(98, 10).Value = 1353.8
$ws.Cells.Item(98, 11).Value = 350.25
$ws.Cells.Item(98, 12).Value = 4061.4
$ws.Cells.Item(98, 13).Value = 1147.75
$ws.Cells.Item(98, 14).Value = -7057.4

$ws.Cells.Item(136, 8).Value = 2038.2
$ws.Cells.Item(136, 9).Value = 1414.5
$ws.Cells.Item(136, 11).Value = 4243.5
$ws.Cells.Item(136, 13).Value = 856.5

$ws.Cells.Item(140, 8).Value = 46154.84
$ws.Cells.Item(140, 9).Value = 56943.75
$ws.Cells.Item(140, 11).Value = 170831.25
$ws.Cells.Item(140, 13).Value = -165651.25


$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 188.10527
$ws.Cells.Item(2, 9).Value = 149.6
$ws.Cells.Item(2, 10).Value = 230.88889
$ws.Cells.Item(2, 11).Value = 149.6
$ws.Cells.Item(2, 12).Value = 230.88889
$ws.Cells.Item(2, 13).Value = -36.59999999999999
$ws.Cells.Item(2, 14).Value = -456.88889

$ws.Cells.Item(132, 8).Value = 3013.0908
$ws.Cells.Item(132, 9).Value = 3044.2222
$ws.Cells.Item(132, 10).Value = 2991.5386
$ws.Cells.Item(132, 11).Value = 9132.6666
$ws.Cells.Item(132, 12).Value = 8974.6158
$ws.Cells.Item(132, 13).Value = -6602.6666
$ws.Cells.Item(132, 14).Value = -14034.6158


$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2327.5715
$ws.Cells.Item(7, 9).Value = 2056.8
$ws.Cells.Item(7, 11).Value = 2056.8
$ws.Cells.Item(7, 13).Value = -1944.8

$ws.Cells.Item(22, 8).Value = 1709.9
$ws.Cells.Item(22, 9).Value = 1585.7142
$ws.Cells.Item(22, 10).Value = 1999.6666
$ws.Cells.Item(22, 11).Value = 1585.7142
$ws.Cells.Item(22, 12).Value = 1999.6666
$ws.Cells.Item(22, 13).Value = -1290.7142
$ws.Cells.Item(22, 14).Value = -2589.6666

$ws.Cells.Item(27, 8).Value = 1709.9
$ws.Cells.Item(27, 9).Value = 1585.7142
$ws.Cells.Item(27, 10).Value = 1999.6666
$ws.Cells.Item(27, 11).Value = 1585.7142
$ws.Cells.Item(27, 12).Value = 1999.6666
$ws.Cells.Item(27, 13).Value = -1478.7142
$ws.Cells.Item(27, 14).Value = -2213.6666

$ws.Cells.Item(40, 8).Value = 3079.8
$ws.Cells.Item(40, 9).Value = 2890
$ws.Cells.Item(40, 11).Value = 2890
$ws.Cells.Item(40, 13).Value = -2754

$ws.Cells.Item(109, 8).Value = 24642.5
$ws.Cells.Item(109, 10).Value = 24642.5
$ws.Cells.Item(109, 12).Value = 24642.5
$ws.Cells.Item(109, 14).Value = -27416.5

$ws.Cells.Item(126, 8).Value = 2327.5715
$ws.Cells.Item(126, 9).Value = 2056.8
$ws.Cells.Item(126, 11).Value = 6170.400000000001
$ws.Cells.Item(126, 13).Value = -3700.400000000001


$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4, 8).Value = 2475
$ws.Cells.Item(4, 9).Value = 2900
$ws.Cells.Item(4, 11).Value = 2900
$ws.Cells.Item(4, 13).Value = -2787

